# Apply "repull data" edits to the dSF column (column F) of Sheet1.
# Source: data/save_data/2023/hatch_thomas.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -3
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = -5
$ws.Range("F15").Value = 4
$ws.Range("F17").Value = -2
$ws.Range("F22").Value = -5
$ws.Range("F24").Value = 13
$ws.Range("F25").Value = -4
